$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F ("想去人数" / want-to-go count)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 648
$wsExhibit.Range("F5").Value = 5037
$wsExhibit.Range("F7").Value = 9509
$wsExhibit.Range("F8").Value = 244
$wsExhibit.Range("F9").Value = 530
$wsExhibit.Range("F10").Value = 88
$wsExhibit.Range("F11").Value = 677
$wsExhibit.Range("F12").Value = 75

# Sheet "全部类型" (All types) - update column F ("想去人数" / want-to-go count)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1214
$wsAll.Range("F3").Value = 648
$wsAll.Range("F7").Value = 5037
$wsAll.Range("F10").Value = 9509
$wsAll.Range("F11").Value = 244
$wsAll.Range("F13").Value = 0
$wsAll.Range("F16").Value = 677
$wsAll.Range("F18").Value = 75
